# Auto commit at 2025-09-27 16:15:13.42
# Refresh the "Metrics" sheet's data values (B2:B13) with the latest
# numbers from the source data feed. Everything downstream ("today" sheet
# formulas that reference Metrics!B2:B13, plus the dependent E/F columns)
# recalculates automatically from these literal updates.

$wb = $excel.ActiveWorkbook

$wsMetrics = $wb.Worksheets.Item("Metrics")

$wsMetrics.Range("B2").Value = 388386.37
$wsMetrics.Range("B3").Value = 313498.31000000006
$wsMetrics.Range("B4").Value = 122788.9
$wsMetrics.Range("B5").Value = 15427
$wsMetrics.Range("B6").Value = 4307637.2499999991
$wsMetrics.Range("B7").Value = 3641025.7899999996
$wsMetrics.Range("B8").Value = 1252154.5800000003
$wsMetrics.Range("B9").Value = 166587
$wsMetrics.Range("B10").Value = 32772961.050999828
$wsMetrics.Range("B11").Value = 19670895.860000007
$wsMetrics.Range("B12").Value = 11533863.469999999
$wsMetrics.Range("B13").Value = 1264214

# Restore the cursor/selection position on the Metrics sheet.
$wsMetrics.Range("F13").Select() | Out-Null

# Move the cursor/selection position on the "today" sheet (the tab that
# was active/selected when the workbook was last saved).
$wsToday = $wb.Worksheets.Item("today")
$wsToday.Activate() | Out-Null
$wsToday.Range("G6").Select() | Out-Null
